# Insert a new weekly price record as row 531 on the single sheet,
# pushing the existing rows 531-602 down to 532-603.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 531 (Excel shifts
# row 531..602 down to 532..603, carrying formatting with it).
$ws.Rows.Item(531).Insert()

# Populate the newly inserted row 531 with the new record.
$ws.Cells.Item(531, 1).Value2  = 3
$ws.Cells.Item(531, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(531, 3).Value2  = "Coquimbo"
$ws.Cells.Item(531, 4).Value2  = 45154
$ws.Cells.Item(531, 5).Value2  = 5
$ws.Cells.Item(531, 6).Value2  = 100112009
$ws.Cells.Item(531, 7).Value2  = "Acelga"
$ws.Cells.Item(531, 8).Value2  = "Sin especificar"
$ws.Cells.Item(531, 9).Value2  = "Primera"
$ws.Cells.Item(531, 10).Value2 = 120
$ws.Cells.Item(531, 11).Value2 = 3000
$ws.Cells.Item(531, 12).Value2 = 3000
$ws.Cells.Item(531, 13).Value2 = 3000
$ws.Cells.Item(531, 14).Value2 = "$/docena de atados (6 kilos)"
$ws.Cells.Item(531, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(531, 16).Value2 = 500
$ws.Cells.Item(531, 17).Value2 = 6
$ws.Cells.Item(531, 18).Value2 = "Hortaliza"
